# Scheduled-runner data refresh: update per-item profit figures (columns
# H:N) on a handful of rows across the ALC, ARM, CRP and WVR sheets.
# Some rows get revised numbers, some previously-blank rows gain a fresh
# set of figures, and a few rows on WVR have their figures cleared back
# out (no longer tracked).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 1379
$ws.Cells.Item(113, 9).Value = 1298.3334
$ws.Cells.Item(113, 11).Value = 1298.3334
$ws.Cells.Item(113, 13).Value = 1955.6666
$ws.Cells.Item(125, 8).Value = 742.5
$ws.Cells.Item(125, 9).Value = 742.5
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 6682.5
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = -4222.5
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(127, 8).Value = 197
$ws.Cells.Item(127, 9).Value = 197
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 11).Value = 591
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 13).Value = 4369
$ws.Cells.Item(128, 8).Value = 0
$ws.Cells.Item(128, 9).Value = 0
$ws.Cells.Item(128, 10).Value = 0
$ws.Cells.Item(128, 11).Value = 0
$ws.Cells.Item(128, 12).Value = 0
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(130, 8).Value = 96653.664
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 96653.664
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 96653.664
$ws.Cells.Item(130, 14).Value = -106693.664
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 4261.1816
$ws.Cells.Item(132, 9).Value = 4261.1816
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 12783.5448
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -10253.5448
$ws.Cells.Item(133, 8).Value = 99995
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 99995
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 99995
$ws.Cells.Item(133, 14).Value = -110115
$ws.Cells.Item(134, 8).Value = 0
$ws.Cells.Item(134, 9).Value = 0
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(135, 8).Value = 902.0833
$ws.Cells.Item(135, 9).Value = 870.4545000000001
$ws.Cells.Item(135, 10).Value = 1250
$ws.Cells.Item(135, 11).Value = 7834.0905
$ws.Cells.Item(135, 12).Value = 11250
$ws.Cells.Item(135, 13).Value = -5299.0905
$ws.Cells.Item(135, 14).Value = -16320
$ws.Cells.Item(136, 8).Value = 99995
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 99995
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 99995
$ws.Cells.Item(136, 14).Value = -110195
$ws.Cells.Item(137, 8).Value = 2220
$ws.Cells.Item(137, 9).Value = 2150
$ws.Cells.Item(137, 10).Value = 2500
$ws.Cells.Item(137, 11).Value = 6450
$ws.Cells.Item(137, 12).Value = 7500
$ws.Cells.Item(137, 13).Value = -3900
$ws.Cells.Item(137, 14).Value = -12600
$ws.Cells.Item(138, 8).Value = 3863.963
$ws.Cells.Item(138, 9).Value = 2250
$ws.Cells.Item(138, 10).Value = 3993.08
$ws.Cells.Item(138, 11).Value = 6750
$ws.Cells.Item(138, 12).Value = 11979.24
$ws.Cells.Item(138, 13).Value = -1610
$ws.Cells.Item(138, 14).Value = -22259.24
$ws.Cells.Item(139, 8).Value = 50780
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 50780
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 50780
$ws.Cells.Item(139, 14).Value = -61060
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 0
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 0
$ws.Cells.Item(141, 8).Value = 14725.23
$ws.Cells.Item(141, 9).Value = 4675.364
$ws.Cells.Item(141, 10).Value = 69999.5
$ws.Cells.Item(141, 11).Value = 14026.092
$ws.Cells.Item(141, 12).Value = 209998.5
$ws.Cells.Item(141, 13).Value = -8846.091999999999
$ws.Cells.Item(141, 14).Value = -220358.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1655.4375
$ws.Cells.Item(32, 9).Value = 1632.5
$ws.Cells.Item(32, 11).Value = 1632.5
$ws.Cells.Item(32, 13).Value = -1345.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(129, 8).Value = 0
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 10).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 12).Value = 0
$ws.Cells.Item(130, 8).Value = 49750
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 49750
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 49750
$ws.Cells.Item(130, 14).Value = -59790
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(132, 8).Value = 2600
$ws.Cells.Item(132, 9).Value = 2677.7778
$ws.Cells.Item(132, 10).Value = 1900
$ws.Cells.Item(132, 11).Value = 8033.3334
$ws.Cells.Item(132, 12).Value = 5700
$ws.Cells.Item(132, 13).Value = -5503.3334
$ws.Cells.Item(132, 14).Value = -10760
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 9).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 11).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(134, 8).Value = 3166.5
$ws.Cells.Item(134, 9).Value = 2899.8
$ws.Cells.Item(134, 10).Value = 4500
$ws.Cells.Item(134, 11).Value = 8699.400000000001
$ws.Cells.Item(134, 12).Value = 13500
$ws.Cells.Item(134, 13).Value = -6164.400000000001
$ws.Cells.Item(134, 14).Value = -18570
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 9).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(137, 8).Value = 0
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(138, 8).Value = 0
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(139, 8).Value = 99994
$ws.Cells.Item(139, 9).Value = 0
$ws.Cells.Item(139, 10).Value = 99994
$ws.Cells.Item(139, 11).Value = 0
$ws.Cells.Item(139, 12).Value = 99994
$ws.Cells.Item(139, 14).Value = -110274
$ws.Cells.Item(140, 8).Value = 99995
$ws.Cells.Item(140, 9).Value = 0
$ws.Cells.Item(140, 10).Value = 99995
$ws.Cells.Item(140, 11).Value = 0
$ws.Cells.Item(140, 12).Value = 99995
$ws.Cells.Item(140, 14).Value = -110355
$ws.Cells.Item(141, 8).Value = 106000
$ws.Cells.Item(141, 9).Value = 50000
$ws.Cells.Item(141, 10).Value = 120000
$ws.Cells.Item(141, 11).Value = 50000
$ws.Cells.Item(141, 12).Value = 120000
$ws.Cells.Item(141, 13).Value = -44820
$ws.Cells.Item(141, 14).Value = -130360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range($ws.Cells.Item(119, 8), $ws.Cells.Item(119, 12)).ClearContents()
$ws.Cells.Item(119, 14).ClearContents()
$ws.Range($ws.Cells.Item(120, 8), $ws.Cells.Item(120, 12)).ClearContents()
$ws.Range($ws.Cells.Item(121, 8), $ws.Cells.Item(121, 12)).ClearContents()
$ws.Range($ws.Cells.Item(122, 8), $ws.Cells.Item(122, 14)).ClearContents()
$ws.Range($ws.Cells.Item(123, 8), $ws.Cells.Item(123, 12)).ClearContents()
$ws.Cells.Item(123, 14).ClearContents()
$ws.Range($ws.Cells.Item(124, 8), $ws.Cells.Item(124, 12)).ClearContents()
$ws.Cells.Item(124, 14).ClearContents()
$ws.Range($ws.Cells.Item(125, 8), $ws.Cells.Item(125, 12)).ClearContents()
$ws.Range($ws.Cells.Item(126, 8), $ws.Cells.Item(126, 14)).ClearContents()
$ws.Range($ws.Cells.Item(127, 8), $ws.Cells.Item(127, 12)).ClearContents()
$ws.Range($ws.Cells.Item(128, 8), $ws.Cells.Item(128, 12)).ClearContents()
$ws.Range($ws.Cells.Item(129, 8), $ws.Cells.Item(129, 12)).ClearContents()
$ws.Cells.Item(129, 14).ClearContents()
$ws.Range($ws.Cells.Item(130, 8), $ws.Cells.Item(130, 12)).ClearContents()
$ws.Range($ws.Cells.Item(131, 8), $ws.Cells.Item(131, 12)).ClearContents()
$ws.Cells.Item(131, 14).ClearContents()
$ws.Range($ws.Cells.Item(132, 8), $ws.Cells.Item(132, 14)).ClearContents()
$ws.Range($ws.Cells.Item(133, 8), $ws.Cells.Item(133, 12)).ClearContents()
$ws.Cells.Item(133, 14).ClearContents()
$ws.Range($ws.Cells.Item(135, 8), $ws.Cells.Item(135, 12)).ClearContents()
$ws.Cells.Item(135, 14).ClearContents()
$ws.Range($ws.Cells.Item(136, 8), $ws.Cells.Item(136, 14)).ClearContents()
$ws.Range($ws.Cells.Item(137, 8), $ws.Cells.Item(137, 12)).ClearContents()
$ws.Range($ws.Cells.Item(138, 8), $ws.Cells.Item(138, 12)).ClearContents()
$ws.Range($ws.Cells.Item(139, 8), $ws.Cells.Item(139, 12)).ClearContents()
$ws.Range($ws.Cells.Item(140, 8), $ws.Cells.Item(140, 12)).ClearContents()
$ws.Range($ws.Cells.Item(141, 8), $ws.Cells.Item(141, 12)).ClearContents()
